# Updated cryptos list on Mon Aug 26 15:56:18 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) and Volume(1h) (E) columns with the latest scrape,
# and reorders the last three coins (Stellar / VeChain / WhiteBITCoin) to
# their new ranking (rows 49-51).
#
# Note: a handful of new Price values are plain decimal numbers (e.g.
# "1.00", "0.998", "38.20"). Assigning those as bare strings via COM lets
# Excel "smart type" them into numeric cells (dropping trailing zeros and
# changing the cell type), same as typing them straight into a worksheet.
# To keep them as text -- matching every other row in this column -- they
# are entered the same way a user would force text in the Excel UI: with a
# leading apostrophe.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, [string]$text) {
    if ($text -match '^[+-]?[0-9]+(\.[0-9]+)?$') {
        $range.Value = "'" + $text
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range("D2")  "63.511.65"
Set-TextValue $ws.Range("E2")  "  -0.83%  "

Set-TextValue $ws.Range("D3")  "2.721.22"

Set-TextValue $ws.Range("E4")  "  -0.05%  "

Set-TextValue $ws.Range("D5")  "560.14"
Set-TextValue $ws.Range("E5")  "  -2.66%  "

Set-TextValue $ws.Range("D6")  "158.21"
Set-TextValue $ws.Range("E6")  "  -0.49%  "

Set-TextValue $ws.Range("E7")  "  -0.02%  "

Set-TextValue $ws.Range("E8")  "  -1.45%  "

Set-TextValue $ws.Range("E9")  "  -2.68%  "

Set-TextValue $ws.Range("E10") "  -0.06%  "

Set-TextValue $ws.Range("D11") "5.64"
Set-TextValue $ws.Range("E11") "  -3.45%  "

Set-TextValue $ws.Range("E12") "  -3.52%  "

Set-TextValue $ws.Range("D13") "3.201.26"
Set-TextValue $ws.Range("E13") "  -1.30%  "

Set-TextValue $ws.Range("D14") "26.56"
Set-TextValue $ws.Range("E14") "  -1.61%  "

Set-TextValue $ws.Range("D15") "63.380.55"
Set-TextValue $ws.Range("E15") "  -0.45%  "

Set-TextValue $ws.Range("E16") "  -3.24%  "

Set-TextValue $ws.Range("D17") "2.722.88"
Set-TextValue $ws.Range("E17") "  -1.38%  "

Set-TextValue $ws.Range("E18") "  +0.29%  "

Set-TextValue $ws.Range("E19") "  -3.96%  "

Set-TextValue $ws.Range("D20") "350.85"
Set-TextValue $ws.Range("E20") "  -1.76%  "

Set-TextValue $ws.Range("E21") "  -4.27%  "

Set-TextValue $ws.Range("E22") "  +0.24%  "

Set-TextValue $ws.Range("E23") "  -3.47%  "

Set-TextValue $ws.Range("D24") "64.46"
Set-TextValue $ws.Range("E24") "  -1.32%  "

Set-TextValue $ws.Range("E25") "  -0.55%  "

Set-TextValue $ws.Range("D26") "1.00"
Set-TextValue $ws.Range("E26") "  +0.03%  "

Set-TextValue $ws.Range("D27") "8.21"
Set-TextValue $ws.Range("E27") "  -4.82%  "

Set-TextValue $ws.Range("E28") "  -2.84%  "

Set-TextValue $ws.Range("E29") "  +8.94%  "

Set-TextValue $ws.Range("E30") "  +0.20%  "

Set-TextValue $ws.Range("E31") "  -2.53%  "

Set-TextValue $ws.Range("D32") "165.45"
Set-TextValue $ws.Range("E32") "  -2.63%  "

Set-TextValue $ws.Range("E33") "  +0.15%  "

Set-TextValue $ws.Range("E34") "  -1.54%  "

Set-TextValue $ws.Range("E35") "  -0.03%  "

Set-TextValue $ws.Range("D36") "4.82"
Set-TextValue $ws.Range("E36") "  -2.74%  "

Set-TextValue $ws.Range("D37") "1.80"
Set-TextValue $ws.Range("E37") "  -0.25%  "

Set-TextValue $ws.Range("D38") "344.66"
Set-TextValue $ws.Range("E38") "  -1.56%  "

Set-TextValue $ws.Range("D39") "0.959"
Set-TextValue $ws.Range("E39") "  -4.33%  "

Set-TextValue $ws.Range("D40") "6.09"
Set-TextValue $ws.Range("E40") "  -3.16%  "

Set-TextValue $ws.Range("E41") "  -3.55%  "

Set-TextValue $ws.Range("D42") "38.20"
Set-TextValue $ws.Range("E42") "  -2.22%  "

Set-TextValue $ws.Range("E43") "  -2.42%  "

Set-TextValue $ws.Range("D44") "20.80"
Set-TextValue $ws.Range("E44") "  -3.40%  "

Set-TextValue $ws.Range("E45") "  -3.03%  "

Set-TextValue $ws.Range("E46") "  -1.01%  "

Set-TextValue $ws.Range("D47") "132.10"
Set-TextValue $ws.Range("E47") "  -3.38%  "

Set-TextValue $ws.Range("D48") "0.998"
Set-TextValue $ws.Range("E48") "  -0.07%  "

# Rows 49-51 reshuffle: WhiteBITCoin / Stellar / VeChain -> Stellar / VeChain / WhiteBITCoin
Set-TextValue $ws.Range("B49") "Stellar"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue $ws.Range("D49") "0.0984"
Set-TextValue $ws.Range("E49") "  -3.34%  "

Set-TextValue $ws.Range("B50") "VeChain"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-TextValue $ws.Range("D50") "0.0245"
Set-TextValue $ws.Range("E50") "  -3.91%  "

Set-TextValue $ws.Range("B51") "WhiteBITCoin"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextValue $ws.Range("D51") "11.04"
Set-TextValue $ws.Range("E51") "  -0.05%  "
